# Updated content in doc2
# Original paragraph 1 text is "Rohini File" followed by a collapsed
# "_GoBack" bookmark. The edit:
#   - wraps "Rohini" in proofErr spellStart/spellEnd markers (as Word's
#     proofing engine would do for a word it doesn't recognize)
#   - splits "Rohini File" into three runs: "Rohini", " File", and a new
#     run ". Update happened here"
#   - keeps the "_GoBack" bookmark collapsed at the end of the paragraph
#     text (after the newly appended sentence)

$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark first so re-creating it further
# along (as part of the XML we insert below) doesn't leave a stray
# duplicate bookmark behind.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate "Rohini File" and replace it (plus re-insert the bookmark after
# the new text) via InsertXML so we get full control over run
# boundaries and the w:proofErr markers, matching exactly what Word's
# OOXML looks like after a spell-checked, multi-run edit.
$target = $d.Content
$target.Find.Execute("Rohini File", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$rng = $target.Duplicate

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Rohini</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> File</w:t></w:r>
            <w:r><w:t>. Update happened here</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$rng.InsertXML($xml)
